# WRI requested web app edits
# Rescales the "Large Output Currency Unit" for the India model from
# crores (10^7) to lakh crores (10^7 * 10^5), updating the related labels.

$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")

# Row 18: label for the large-unit scale changes from "2018 crores" to
# "2018 lakh crores".
$about.Range("A18").Value = "2018 lakh crores"

# Row 19: the scale factor formula grows by another factor of 10^5 (lakh),
# and its units label changes to match.
$about.Range("A19").Formula = "=10^7*10^5"
$about.Range("B19").Value = "rupees per lakh crore"

# Update the active selection shown on the About sheet to A20 (matches the
# saved view after the edit) and clear the scrolled "top-left cell".
$about.Activate() | Out-Null
$about.Range("A20").Select() | Out-Null
